$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IR")

# Set column F width to match updated diff value (14.3 characters).
# Excel's ColumnWidth setter snaps to whole-pixel boundaries (character
# width is stored internally as pixels/MaxDigitWidth), so 13.6 is the
# input that lands on the pixel boundary closest to the target 14.3.
$ws.Range("F1").ColumnWidth = 13.6

# Row 4: Non Cash Items (Other)
$ws.Range("B4").Value = 238100000.0
$ws.Range("F4").Value = 10300000.0

# Row 6: Change in inventories
$ws.Range("B6").Value = 138300000.0
$ws.Range("C6").Value = 170800000.0
$ws.Range("D6").Value = 139700000.0
$ws.Range("E6").Value = 45900000.0
$ws.Range("F6").Value = 28600000.0

# Row 7: Accounts Payable Change
$ws.Range("B7").Value = -63000000.0
$ws.Range("F7").Value = 63100000.0

# Row 8: Change in payables and accrued liability
$ws.Range("B8").Value = 2654000000.0
$ws.Range("C8").Value = 2745000000.0
$ws.Range("D8").Value = 2058800000.0
$ws.Range("E8").Value = 1423800000.0
$ws.Range("F8").Value = 723400000.0

# Row 11: Capital expenditures
$ws.Range("B11").Value = -44800000.0
$ws.Range("F11").Value = -35600000.0

# Row 23: Stock Based Compensation
$ws.Range("B23").Value = 71200000.0
$ws.Range("F23").Value = 15100000.0

# Row 25: Assets Liabilities Change (Total)
$ws.Range("B25").Value = 269900000.0
$ws.Range("F25").Value = 31200000.0
